$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells keep their original text format so that
# numeric-looking strings (e.g. "1.010", "12.90") are not coerced to numbers.
$targetCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "E31", "E32", "D33", "E33", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "E49", "D50", "E50", "B51", "C51", "D51", "E51")
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.596.72"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.116.04"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "336.48"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "0.5246"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "0.4551"
$ws.Range("E8").Value = "  +3.23%  "
$ws.Range("D9").Value = "54.61"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "0.09137"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "24.47"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "2.115.73"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "6.848"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "8.146"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").Value = "0.00001178"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "97.15"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "1.011"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "19.41"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "6.305"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "30.656.28"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "12.90"
$ws.Range("E24").Value = "  +5.12%  "
$ws.Range("D25").Value = "2.351"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "2.373.22"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "164.39"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "2.548"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "134.63"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "1.648"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("D35").Value = "3.945"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "10.63"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("D37").Value = "5.868"
$ws.Range("E37").Value = "  +7.43%  "
$ws.Range("D38").Value = "0.02625"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "0.06837"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "0.2325"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("D41").Value = "12.57"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "0.6888"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("D45").Value = "0.6483"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "2.309"
$ws.Range("E46").Value = "  +5.43%  "
$ws.Range("D47").Value = "0.00000000369"
$ws.Range("E47").Value = "  +22.07%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "83.38"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07215"
$ws.Range("E51").Value = "  +2.51%  "
